$wb = $excel.ActiveWorkbook

# Sheet: IBT
$wsIBT = $wb.Worksheets.Item("IBT")
$wsIBT.Range("C3").Value = 1
$wsIBT.Range("C4").Value = 1
$wsIBT.Range("C6").Value = 1
$wsIBT.Range("C7").Value = 1
[void]$wsIBT.Range("C1").Select()

# Sheet: Pension
$wsPension = $wb.Worksheets.Item("Pension")
$wsPension.Range("C3").Value = 1
$wsPension.Range("C4").Value = 1
$wsPension.Range("C6").Value = 1
$wsPension.Range("C7").Value = 1
[void]$wsPension.Range("C2").Select()

# Sheet: Retirement
$wsRetirement = $wb.Worksheets.Item("Retirement")
$wsRetirement.Range("C3").Value = 1
$wsRetirement.Range("C4").Value = 1
$wsRetirement.Range("C6").Value = 1
$wsRetirement.Range("C7").Value = 1
[void]$wsRetirement.Range("C2").Select()

# Restore IBT as the active sheet with its selection active, matching tabSelected="1"
[void]$wsIBT.Select()
[void]$wsIBT.Range("C1").Select()
